# Adds a new "2022-Q4" worksheet (with its fund-holding detail data) right
# after the "总计" sheet, and inserts the corresponding 2022-Q4 summary row
# into the "总计" sheet, pushing the older quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Part 1: Insert a new "2022-Q4" worksheet right after "总计"
# (i.e. right before the current "2022-Q1" sheet)
# ---------------------------------------------------------------
$wsQ1Ref = $wb.Worksheets.Item(2)
$newWs = $wb.Worksheets.Add($wsQ1Ref)
$newWs.Name = "2022-Q4"

# Re-fetch the "2022-Q1" sheet fresh by name - the reference used above as
# the insertion anchor can go stale once the new sheet has been inserted.
$wsQ1 = $wb.Worksheets.Item("2022-Q1")

# Copy the header style (bold + border) from the template sheet onto row 1
# (B1:H1) and onto the index column (A2:A3) of the new sheet.
$wsQ1.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$wsQ1.Range("A2:A3").Copy()
$newWs.Range("A2:A3").PasteSpecial(-4122)

# Header row text
$newWs.Cells.Item(1,2).Value = "基金代码"
$newWs.Cells.Item(1,3).Value = "基金名称"
$newWs.Cells.Item(1,4).Value = "基金规模"
$newWs.Cells.Item(1,5).Value = "股票总仓位"
$newWs.Cells.Item(1,6).Value = "仓位占比"
$newWs.Cells.Item(1,7).Value = "持有市值(亿元)"
$newWs.Cells.Item(1,8).Value = "仓位排名"

# Index column (A2/A3)
$newWs.Cells.Item(2,1).Value = 0
$newWs.Cells.Item(3,1).Value = 1

# Text-typed data cells - force text storage by briefly switching the
# number format to "@" before assigning the value (otherwise values such
# as "005618" or "0.13" get silently converted to numbers), then drop the
# leftover number-format style so the cell keeps the default, unstyled look.
$textCells = @(
  @{r=2; c=2; v="005618"},
  @{r=2; c=3; v="融通红利机会主题精选灵活配置混合A"},
  @{r=2; c=4; v="0.13"},
  @{r=2; c=5; v="84.29"},
  @{r=2; c=6; v="2.93"},
  @{r=2; c=7; v="0.0038"},
  @{r=3; c=2; v="005619"},
  @{r=3; c=3; v="融通红利机会主题精选灵活配置混合C"},
  @{r=3; c=4; v="0.02"},
  @{r=3; c=5; v="84.29"},
  @{r=3; c=6; v="2.93"},
  @{r=3; c=7; v="0.0006"}
)
foreach ($tc in $textCells) {
  $cc = $newWs.Cells.Item($tc.r, $tc.c)
  $cc.NumberFormat = "@"
  $cc.Value = $tc.v
  $cc.ClearFormats()
}

# Numeric rank column
$newWs.Cells.Item(2,8).Value = 8
$newWs.Cells.Item(3,8).Value = 8

# ---------------------------------------------------------------
# Part 2: Update the "总计" (summary) sheet - insert a new row for
# 2022-Q4 right after the header row, pushing all other quarters down.
# ---------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Give A2 the same style as the other index-column cells.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 2
$wsTotal.Cells.Item(2,4).Value = 0

# Re-sequence the index column for the rows that shifted down.
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(6,1).Value = 4
$wsTotal.Cells.Item(7,1).Value = 5
$wsTotal.Cells.Item(8,1).Value = 6

Write-Output "edit complete"
